# Add team record columns (Wins, Losses, Ties) to the data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting from an existing header cell (A1) onto the
# three new header cells so they match the rest of the header row
# (bold font, borders, centered alignment).
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Set the new header labels.
$ws.Range("AD1").Value2 = "Wins"
$ws.Range("AE1").Value2 = "Losses"
$ws.Range("AF1").Value2 = "Ties"

# Fill in the team record (Wins=72, Losses=90, Ties=0) for every data row
# (rows 2 through 53).
$lastRow = 53
$rowCount = $lastRow - 1
$data = New-Object 'object[,]' $rowCount,3
for ($i = 0; $i -lt $rowCount; $i++) {
    $data[$i,0] = 72
    $data[$i,1] = 90
    $data[$i,2] = 0
}
$ws.Range("AD2:AF53").Value2 = $data
